# Revert "Adding CDI benchmark in the summary tab, related to the closed
# operations results" — restore the pre-filter, pre-summary-row state of
# the "Extrato" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear the autofilter criteria (CIEL3 / Acoes) but keep the filter
#    button range itself (A1:P35 once row 36 is gone).
# ---------------------------------------------------------------------
$ws.ShowAllData()

# ---------------------------------------------------------------------
# 2. Un-hide all the data rows that the closed-operations filter had
#    hidden (every row except the header and row 8, which was already
#    visible).
# ---------------------------------------------------------------------
$hiddenRows = @(2,3,4,5,6,7,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $false
}

# ---------------------------------------------------------------------
# 3. Restore the original transaction data that the CDI-benchmark
#    commit had overwritten with placeholder numbers.
# ---------------------------------------------------------------------

# Row 3
$ws.Range("A3").Value = 41964
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 10
$ws.Range("K3").Value = 10

# Row 7
$ws.Range("A7").Value = 42107
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 15
$ws.Range("K7").Value = 10

# Row 8
$ws.Range("A8").Value = 42111
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 15
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 50

# Row 21
$ws.Range("A21").Value = 44312

# Row 22
$ws.Range("A22").Value = 44312

# Row 23
$ws.Range("A23").Value = 44312

# Row 34
$ws.Range("A34").Value = 42695
$ws.Range("H34").Value = 100
$ws.Range("I34").Value = 5
$ws.Range("K34").Value = 10

# Row 35
$ws.Range("A35").Value = 42473
$ws.Range("H35").Value = 100
$ws.Range("I35").Value = 20
$ws.Range("K35").Value = 10

# ---------------------------------------------------------------------
# 4. Remove row 36 entirely (it was a duplicated/extra summary row
#    added by the reverted commit).
# ---------------------------------------------------------------------
$ws.Rows.Item(36).Delete()

# ---------------------------------------------------------------------
# 5. Restore the view: selection on the bottom-right frozen pane sits
#    on A35 (last data row) instead of the default A1.
# ---------------------------------------------------------------------
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A35").Select()
